$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new blank column at N ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 9.86

# Select J16 on this sheet and make it the active/selected sheet/tab
$ws.Activate()
$ws.Range("J16").Select()

# --- "Transactions" sheet: it is no longer the tab-selected sheet ---
# (selecting "Repayment schedule" above already moves the active tab away
# from "Transactions")
